# Update the "Förändrad" (last-changed) date column (C) for rows 2-97
# from 2023-09-17 (serial 45186) to 2023-09-19 (serial 45188).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C97").Value = 45188
